$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row -> [C_new, D_new, E_new]
$updates = @{
    2   = @(766341, 155968, 1429248811)
    13  = @(187861, 33264, 1168232002)
    19  = @(27514, 4210, 132331474)
    21  = @(175242, 38057, 316822468)
    41  = @(126946, 24257, 662711335)
    57  = @(31597, 6018, 162625330)
    81  = @(88357, 16598, 499681283)
    88  = @(71278, 12436, 110315721)
    121 = @(1306396, 220388, 2275480279)
    129 = @(633775, 104970, 3435036140)
    132 = @(586022, 90786, 3472775626)
    151 = @(39935, 7155, 60390275)
    156 = @(12412, 2144, 40642757)
    171 = @(95828, 18309, 490704985)
    178 = @(515888, 115380, 891213559)
    186 = @(236840, 46319, 1190035149)
    237 = @(283328, 49797, 1438555325)
    240 = @(205932, 33984, 1069916235)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}
